$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B (descripcion) and E (parametros) values between row 2 and row 3,
# leaving A, C, D untouched.
$b2 = $ws.Range("B2").Value2
$b3 = $ws.Range("B3").Value2
$e2 = $ws.Range("E2").Value2
$e3 = $ws.Range("E3").Value2

$ws.Range("B2").Value = $b3
$ws.Range("B3").Value = $b2
$ws.Range("E2").Value = $e3
$ws.Range("E3").Value = $e2

# Update the sheet view: clear the frozen/scrolled top-left cell and move
# the active selection to G7.
$ws.Range("G7").Select()
